# The "_GoBack" bookmark (Word's last-edit-position marker) currently sits
# at the end of the "Week: 12" paragraph. Move it to the very start of the
# first paragraph ("Fluency Review Template"), right before its run content.
#
# Note: "_GoBack" is a hidden/special bookmark that Word's Bookmarks
# collection excludes from enumeration (Count/ForEach), but it can still be
# looked up by name directly.

$d = $word.ActiveDocument

# Remove the bookmark from its current location.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Collapse a range to the very beginning of the first paragraph.
$firstPara = $d.Paragraphs(1)
$startRange = $firstPara.Range.Duplicate
$startRange.Collapse(1)

# Re-create the bookmark there via a WordOpenXML fragment (Bookmarks.Add
# mishandles collapsed ranges at document position 0, fusing bookmarkStart
# and bookmarkEnd across the following paragraph instead of keeping them
# adjacent), so splice the pair in directly as WordprocessingML.
$bookmarkFragment = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="99" w:name="_GoBack"/><w:bookmarkEnd w:id="99"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$startRange.InsertXML($bookmarkFragment)
